$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. INVENTARIO sheet: rebuild as a product-inventory table (instead of the
#    old blank CLIENTE-insert-statement scratch sheet).
# ---------------------------------------------------------------------------
$inv = $wb.Worksheets.Item("INVENTARIO")

# Wipe the previous scratch content (headers in A1:H1 + helper formulas I2:I6)
$inv.Cells.Clear() | Out-Null

# Header row
$inv.Range("A1").Value = "fecha_caducidad"
$inv.Range("B1").Value = "cantidades_disponibles"
$inv.Range("C1").Value = "medida"
$inv.Range("D1").Value = "unidad_medida"
$inv.Range("E1").Value = "especialidad"
$inv.Range("F1").Value = "id_sede"
$inv.Range("G1").Value = "nombre_prod"

# Row 2
$inv.Range("A2").Value = " '2020-01-01'"
$inv.Range("B2").Value = 20
$inv.Range("C2").Value = 1
$inv.Range("D2").Value = " 'libra'"
$inv.Range("E2").Value = " 'grano'"
$inv.Range("F2").Value = 1
$inv.Range("G2").Value = " 'pasta larga'"

# Row 3
$inv.Range("A3").Value = " '2021-03-07'"
$inv.Range("B3").Value = 30
$inv.Range("C3").Value = 500
$inv.Range("D3").Value = " 'ml'"
$inv.Range("E3").Value = " 'bebida'"
$inv.Range("F3").Value = 1
$inv.Range("G3").Value = " 'vino'"

# Row 4
$inv.Range("A4").Value = " '2021-04-09'"
$inv.Range("B4").Value = 10
$inv.Range("C4").Value = 100
$inv.Range("D4").Value = " 'grs'"
$inv.Range("E4").Value = " 'condimento'"
$inv.Range("F4").Value = 2
$inv.Range("G4").Value = " 'pimienta'"

# Row 5
$inv.Range("A5").Value = " '2020-12-11'"
$inv.Range("B5").Value = 5
$inv.Range("C5").Value = 1
$inv.Range("D5").Value = " 'libra'"
$inv.Range("E5").Value = " 'carne'"
$inv.Range("F5").Value = 3
$inv.Range("G5").Value = " 'lomo de cerdo'"

# Formula column (H): row 2 stands alone, rows 3-6 form a shared formula
# group (row 6 has no source data, producing an all-empty INSERT).
$inv.Range("H2").Formula = '=+"INSERT INTO cliente("&$A$1&", "&$B$1&", "&$C$1&", "&$D$1&", "&$E$1&", "&$F$1&", "&$G$1&") VALUES ("&A2&","&B2&","&C2&","&D2&","&E2&","&F2&","&G2&");"'
$inv.Range("H3:H6").Formula = '=+"INSERT INTO cliente("&$A$1&", "&$B$1&", "&$C$1&", "&$D$1&", "&$E$1&", "&$F$1&", "&$G$1&") VALUES ("&A3&","&B3&","&C3&","&D3&","&E3&","&F3&","&G3&");"'

$inv.Range("H2:H5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. MENU sheet: the shared-formula range on G6 over-reported its extent
#    (G3:G6 instead of just G6) - tighten it back to the single cell.
# ---------------------------------------------------------------------------
$menu = $wb.Worksheets.Item("MENU")
$menu.Range("G6").Formula = '=+"INSERT INTO cliente("&$A$1&", "&$B$1&", "&$C$1&", "&$D$1&", "&$E$1&", "&$F$1&") VALUES ("&A6&","&B6&","&C6&","&D6&","&E6&","&F6&")"'

# ---------------------------------------------------------------------------
# 3. Active-sheet / selection bookkeeping: the workbook was left with
#    PRODUCTOS_POR_FACTURA selected at E13:E20 - move that selection to
#    D2 and make FACTURACION (selection E3) the active tab instead.
# ---------------------------------------------------------------------------
$prod = $wb.Worksheets.Item("PRODUCTOS_POR_FACTURA")
$prod.Range("D2").Select() | Out-Null

$fact = $wb.Worksheets.Item("FACTURACION")
$fact.Range("E3").Select() | Out-Null
